# Validate_Dulicates.xlsx - "Trials to avoide compilation error"
#
# The birth-year column (D) on Sheet1 is normalised so every row reports
# the same year (1990), the active selection is moved to G6, and the
# sheet's page setup is switched to portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# bYear column (D2:D6) -> 1990 for every data row
$ws.Range("D2").Value = 1990
$ws.Range("D3").Value = 1990
$ws.Range("D4").Value = 1990
$ws.Range("D5").Value = 1990
$ws.Range("D6").Value = 1990

# Move the active selection from G1 to G6
$ws.Range("G6").Select() | Out-Null

# Force portrait page orientation (xlPortrait = 1)
$ws.PageSetup.Orientation = 1
